$excel.Width = 23260
$excel.Height = 12580
$excel.Left = 14360
$excel.Top = 2700
